$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1 ("The data set has been collected from the Bureau of
# Meteorology ... " paragraph):
#   - The hyperlink text "Bureau of Me" + "t" + "eorology" (three runs)
#     becomes a single run "Bureau of Meteorology".
#   - The word "can" becomes "and" (" can contains" -> " and contains"),
#     and that word ends up as its own run, flanked by a run containing
#     just a space and a run with the remaining sentence.
# -----------------------------------------------------------------------

# Step 1: Re-apply the same visible text to the hyperlink run so the three
# runs that build up "Bureau of Meteorology" get consolidated into one.
# (Find/Replace cannot cross the hyperlink's own run/field boundary when a
# replacement is requested together with trailing text, so this is kept
# to text that lives fully inside the hyperlink.)
$r1 = $d.Content.Find.Execute("Bureau of Meteorology", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "Bureau of Meteorology", 2)

# Step 2: Locate the word "can" (whole word) that follows the hyperlink and
# swap it for "and".
$rng = $d.Content
$found = $rng.Find.Execute("can", $true, $true, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "and"
    # Nudge formatting on the freshly written word so the run that holds it
    # stays distinct from its neighbours instead of being re-coalesced with
    # them (mirrors how the original edit left "and" as its own run).
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}

# -----------------------------------------------------------------------
# Change 2 (Technologies / Visualisation line): the run "Visualisation:
# JavaScript (D3, " + "Plotly" (wrapped in proofErr spell-check tags) +
# ", Leaflet) HTML/CSS: Bootstrap" collapses into a single plain run with
# no proofErr markers. The visible text is unchanged, so retyping it via
# Find/Replace consolidates the runs and drops the spell-check tags.
# -----------------------------------------------------------------------
$d.Content.Find.Execute("Visualisation: JavaScript (D3, Plotly, Leaflet) HTML/CSS: Bootstrap", `
                         $true, $false, $false, $false, $false, $true, 1, $false, `
                         "Visualisation: JavaScript (D3, Plotly, Leaflet) HTML/CSS: Bootstrap", 2)
